{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Inserts three new bullet paragraphs into the \"PARTNER - Siege Analytics\"\n// section, immediately after the \"Research & Data Analytics Leadership\"\n// paragraph and before the \"Conceived, architected, engineered...\" bullet.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the anchor paragraph by its exact text.\nconst anchorText = \"Research & Data Analytics Leadership\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph: \" + anchorText);\n}\n\n// New bullet lines to add, in document order.\nconst newLines = [\n  \"\\u2022 Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters\",\n  \"\\u2022 Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States\",\n  \"\\u2022 Algorithm reduced mapping costs by 75%, saving campaigns and organizations $5M+ and enabling smaller nonprofits to conduct redistricting analysis\"\n];\n\n// Insert each line after the anchor, keeping insertion order by chaining\n// \"After\" inserts onto the previously inserted paragraph.\nlet insertAfter = anchor;\nfor (const line of newLines) {\n  insertAfter = insertAfter.insertParagraph(line, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word / $d (ActiveDocument) are pre-seeded by the host.\n#\n# Inserts three new bullet paragraphs into the \"PARTNER - Siege Analytics\"\n# section, immediately after the \"Research & Data Analytics Leadership\"\n# paragraph and before the \"Conceived, architected, engineered...\" bullet.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph by its exact (trimmed) text.\n$anchorText = \"Research & Data Analytics Leadership\"\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $anchorText) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find anchor paragraph: $anchorText\"\n}\n\n# New bullet lines to add, in document order. Word's bullet glyph (U+2022)\n# isn't reachable via a literal in this shell, so build it with [char].\n$bullet = [char]0x2022\n$lines = @(\n    \"$bullet Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters\",\n    \"$bullet Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States\",\n    \"$bullet Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis\"\n)\n\n# Insert each line after the anchor, advancing to the freshly created\n# paragraph each time so the three lines land in order.\n$cur = $target\nforeach ($line in $lines) {\n    $r = $cur.Range\n    $r.Collapse(0)\n    $r.InsertParagraphAfter()\n    $cur = $cur.Next()\n    $cur.Range.InsertAfter($line)\n}\n"}
